# "Generate Report for Archive"
#
# Localization status moved on from handoff: the shared "Ready for handoff"
# status text becomes "In Translation" everywhere it is used (Overview sheet's
# per-language status columns, and each language sheet's Status column), and
# those now-narrower status columns are resized to match.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E & F), rows 2-3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = $newStatus

# --- Per-language sheets: Status column (C), rows 2-3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = $newStatus

# --- Narrow the affected columns to fit the shorter status text ---
# (target stored width ~= 13.41 "characters"; ColumnWidth = 12.5 rounds to the
# closest representable stored width for this column-width model)
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
